# PowerShell COM-interop script to update cryptos.xlsx price (D) and
# volume (E) columns, mirroring the GitHub Actions "Updated cryptos list"
# commit. All values are plain display text in the source sheet (t="inlineStr"),
# so a leading apostrophe is used to force Excel to keep them as text instead
# of auto-converting numeric-looking strings (e.g. "596.91") into numbers,
# which would otherwise drop meaningful trailing zeros (e.g. "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "'68.612.70"
$ws.Range("E2").Value = "'  +1.54%  "

# Row 3
$ws.Range("D3").Value = "'3.782.37"
$ws.Range("E3").Value = "'  +0.56%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'596.91"

# Row 6
$ws.Range("D6").Value = "'169.76"
$ws.Range("E6").Value = "'  -0.83%  "

# Row 7
$ws.Range("D7").Value = "'3.781.19"
$ws.Range("E7").Value = "'  +0.55%  "

# Row 8
$ws.Range("E8").Value = "'  -0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "'  -0.54%  "

# Row 10
$ws.Range("E10").Value = "'  -1.49%  "

# Row 11
$ws.Range("D11").Value = "'6.54"
$ws.Range("E11").Value = "'  +0.86%  "

# Row 12
$ws.Range("E12").Value = "'  -0.60%  "

# Row 13
$ws.Range("E13").Value = "'  -2.69%  "

# Row 14
$ws.Range("D14").Value = "'36.89"
$ws.Range("E14").Value = "'  +0.18%  "

# Row 15
$ws.Range("D15").Value = "'4.417.17"
$ws.Range("E15").Value = "'  +0.54%  "

# Row 16
$ws.Range("D16").Value = "'3.768.35"
$ws.Range("E16").Value = "'  +0.24%  "

# Row 17
$ws.Range("D17").Value = "'68.578.16"
$ws.Range("E17").Value = "'  +1.38%  "

# Row 18
$ws.Range("D18").Value = "'18.19"
$ws.Range("E18").Value = "'  -3.88%  "

# Row 19
$ws.Range("E19").Value = "'  -2.49%  "

# Row 20
$ws.Range("E20").Value = "'  -0.11%  "

# Row 21
$ws.Range("D21").Value = "'10.98"
$ws.Range("E21").Value = "'  +3.97%  "

# Row 22
$ws.Range("D22").Value = "'469.25"
$ws.Range("E22").Value = "'  -0.17%  "

# Row 23
$ws.Range("D23").Value = "'0.705"
$ws.Range("E23").Value = "'  -2.36%  "

# Row 24
$ws.Range("D24").Value = "'84.79"
$ws.Range("E24").Value = "'  +1.12%  "

# Row 25
$ws.Range("E25").Value = "'  -2.98%  "

# Row 26
$ws.Range("E26").Value = "'  +0.47%  "

# Row 27
$ws.Range("D27").Value = "'12.24"
$ws.Range("E27").Value = "'  +0.14%  "

# Row 28
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "'  -1.35%  "

# Row 29
$ws.Range("E29").Value = "'  +0.05%  "

# Row 30
$ws.Range("D30").Value = "'3.928.60"
$ws.Range("E30").Value = "'  +0.47%  "

# Row 31
$ws.Range("E31").Value = "'  -3.37%  "

# Row 32
$ws.Range("D32").Value = "'7.43"
$ws.Range("E32").Value = "'  -3.33%  "

# Row 33
$ws.Range("E33").Value = "'  -1.02%  "

# Row 34
$ws.Range("D34").Value = "'30.16"
$ws.Range("E34").Value = "'  -0.96%  "

# Row 35
$ws.Range("D35").Value = "'9.38"
$ws.Range("E35").Value = "'  +2.46%  "

# Row 37
$ws.Range("D37").Value = "'3.736.40"
$ws.Range("E37").Value = "'  +0.23%  "

# Row 38
$ws.Range("E38").Value = "'  -2.73%  "

# Row 39
$ws.Range("E39").Value = "'  -9.33%  "

# Row 40
$ws.Range("E40").Value = "'  +1.17%  "

# Row 41
$ws.Range("D41").Value = "'1.00"

# Row 42
$ws.Range("D42").Value = "'5.85"
$ws.Range("E42").Value = "'  -1.03%  "

# Row 43
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "'  -0.15%  "

# Row 44
$ws.Range("D44").Value = "'0.311"
$ws.Range("E44").Value = "'  -0.84%  "

# Row 45
$ws.Range("E45").Value = "'  +0.01%  "

# Row 46
$ws.Range("D46").Value = "'1.97"
$ws.Range("E46").Value = "'  +0.47%  "

# Row 47
$ws.Range("D47").Value = "'8.61"
$ws.Range("E47").Value = "'  -1.21%  "

# Row 48
$ws.Range("D48").Value = "'42.81"
$ws.Range("E48").Value = "'  +9.88%  "

# Row 49
$ws.Range("D49").Value = "'401.55"
$ws.Range("E49").Value = "'  +0.85%  "

# Row 50
$ws.Range("D50").Value = "'45.81"
$ws.Range("E50").Value = "'  +0.07%  "

# Row 51
$ws.Range("D51").Value = "'145.30"
$ws.Range("E51").Value = "'  +2.67%  "
